# Update recomputed time-allocation probabilities in the South Carolina_A
# team-specific transition matrix (row-stochastic values per starting state).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Starting_State = 0)
$ws.Cells.Item(2, 2).Value = 0.1790322580645161  # B2
$ws.Cells.Item(2, 3).Value = 0.5596774193548387  # C2
$ws.Cells.Item(2, 10).Value = 0.01935483870967742  # J2
$ws.Cells.Item(2, 16).Value = 0.1290322580645161  # P2
$ws.Cells.Item(2, 19).Value = 0.1129032258064516  # S2

# Row 3 (Starting_State = 1)
$ws.Cells.Item(3, 2).Value = 0.01966292134831461  # B3
$ws.Cells.Item(3, 3).Value = 0.01404494382022472  # C3
$ws.Cells.Item(3, 10).Value = 0.05898876404494382  # J3
$ws.Cells.Item(3, 16).Value = 0.7219101123595506  # P3
$ws.Cells.Item(3, 19).Value = 0.1853932584269663  # S3

# Row 4 (Starting_State = 2)
$ws.Cells.Item(4, 10).Value = 0.08695652173913043  # J4
$ws.Cells.Item(4, 16).Value = 0.6811594202898551  # P4
$ws.Cells.Item(4, 19).Value = 0.2318840579710145  # S4

# Row 5 (Starting_State = 3)
$ws.Cells.Item(5, 16).Value = 1  # P5

# Row 6 (Starting_State = 4)
$ws.Cells.Item(6, 2).Value = 0.05429864253393665  # B6
$ws.Cells.Item(6, 4).Value = 0.002262443438914027  # D6
$ws.Cells.Item(6, 6).Value = 0.08597285067873303  # F6
$ws.Cells.Item(6, 10).Value = 0.2760180995475113  # J6
$ws.Cells.Item(6, 15).Value = 0.01357466063348416  # O6
$ws.Cells.Item(6, 17).Value = 0.1900452488687783  # Q6
$ws.Cells.Item(6, 18).Value = 0.06561085972850679  # R6
$ws.Cells.Item(6, 19).Value = 0.3122171945701357  # S6

# Row 7 (Starting_State = 5)
$ws.Cells.Item(7, 2).Value = 0.1137440758293839  # B7
$ws.Cells.Item(7, 4).Value = 0.01895734597156398  # D7
$ws.Cells.Item(7, 6).Value = 0.04265402843601896  # F7
$ws.Cells.Item(7, 10).Value = 0.1587677725118483  # J7
$ws.Cells.Item(7, 15).Value = 0.02369668246445497  # O7
$ws.Cells.Item(7, 17).Value = 0.1658767772511848  # Q7
$ws.Cells.Item(7, 18).Value = 0.08767772511848342  # R7
$ws.Cells.Item(7, 19).Value = 0.3886255924170616  # S7

# Row 8 (Starting_State = 6)
$ws.Cells.Item(8, 2).Value = 0.1172839506172839  # B8
$ws.Cells.Item(8, 4).Value = 0.01604938271604938  # D8
$ws.Cells.Item(8, 6).Value = 0.05802469135802469  # F8
$ws.Cells.Item(8, 10).Value = 0.1395061728395062  # J8
$ws.Cells.Item(8, 15).Value = 0.02098765432098765  # O8
$ws.Cells.Item(8, 17).Value = 0.154320987654321  # Q8
$ws.Cells.Item(8, 18).Value = 0.08518518518518518  # R8
$ws.Cells.Item(8, 19).Value = 0.408641975308642  # S8

# Row 9 (Starting_State = 7)
$ws.Cells.Item(9, 2).Value = 0.1337386018237082  # B9
$ws.Cells.Item(9, 4).Value = 0.00911854103343465  # D9
$ws.Cells.Item(9, 6).Value = 0.05167173252279635  # F9
$ws.Cells.Item(9, 10).Value = 0.1458966565349544  # J9
$ws.Cells.Item(9, 15).Value = 0.00911854103343465  # O9
$ws.Cells.Item(9, 17).Value = 0.1793313069908815  # Q9
$ws.Cells.Item(9, 18).Value = 0.09422492401215805  # R9
$ws.Cells.Item(9, 19).Value = 0.3768996960486322  # S9

# Row 10 (Starting_State = 8)
$ws.Cells.Item(10, 2).Value = 0.1080974842767296  # B10
$ws.Cells.Item(10, 4).Value = 0.01886792452830189  # D10
$ws.Cells.Item(10, 5).Value = 0.0003930817610062893  # E10
$ws.Cells.Item(10, 6).Value = 0.06525157232704402  # F10
$ws.Cells.Item(10, 10).Value = 0.1525157232704402  # J10
$ws.Cells.Item(10, 15).Value = 0.02830188679245283  # O10
$ws.Cells.Item(10, 17).Value = 0.2028301886792453  # Q10
$ws.Cells.Item(10, 18).Value = 0.07861635220125786  # R10
$ws.Cells.Item(10, 19).Value = 0.345125786163522  # S10

# Row 11 (Starting_State = 9)
$ws.Cells.Item(11, 7).Value = 0.1595092024539877  # G11
$ws.Cells.Item(11, 10).Value = 0.08435582822085889  # J11
$ws.Cells.Item(11, 11).Value = 0.2070552147239264  # K11
$ws.Cells.Item(11, 12).Value = 0.5383435582822086  # L11
$ws.Cells.Item(11, 19).Value = 0.01073619631901841  # S11

# Row 12 (Starting_State = 10)
$ws.Cells.Item(12, 7).Value = 0.7438692098092643  # G12
$ws.Cells.Item(12, 10).Value = 0.1907356948228883  # J12
$ws.Cells.Item(12, 11).Value = 0.0108991825613079  # K12
$ws.Cells.Item(12, 12).Value = 0.02997275204359673  # L12
$ws.Cells.Item(12, 19).Value = 0.02452316076294278  # S12

# Row 13 (Starting_State = 11)
$ws.Cells.Item(13, 6).Value = 0.0136986301369863  # F13
$ws.Cells.Item(13, 7).Value = 0.7397260273972602  # G13
$ws.Cells.Item(13, 10).Value = 0.2054794520547945  # J13
$ws.Cells.Item(13, 19).Value = 0.0410958904109589  # S13

# Row 15 (Starting_State = 13)
$ws.Cells.Item(15, 6).Value = 0.01345291479820628  # F15
$ws.Cells.Item(15, 8).Value = 0.1704035874439462  # H15
$ws.Cells.Item(15, 9).Value = 0.05829596412556054  # I15
$ws.Cells.Item(15, 10).Value = 0.352017937219731  # J15
$ws.Cells.Item(15, 11).Value = 0.08071748878923767  # K15
$ws.Cells.Item(15, 13).Value = 0.008968609865470852  # M15
$ws.Cells.Item(15, 15).Value = 0.06278026905829596  # O15
$ws.Cells.Item(15, 19).Value = 0.2533632286995516  # S15

# Row 16 (Starting_State = 14)
$ws.Cells.Item(16, 6).Value = 0.01912568306010929  # F16
$ws.Cells.Item(16, 8).Value = 0.1693989071038251  # H16
$ws.Cells.Item(16, 9).Value = 0.06010928961748634  # I16
$ws.Cells.Item(16, 10).Value = 0.4562841530054645  # J16
$ws.Cells.Item(16, 11).Value = 0.1092896174863388  # K16
$ws.Cells.Item(16, 13).Value = 0.02185792349726776  # M16
$ws.Cells.Item(16, 15).Value = 0.04098360655737705  # O16
$ws.Cells.Item(16, 19).Value = 0.1229508196721311  # S16

# Row 17 (Starting_State = 15)
$ws.Cells.Item(17, 6).Value = 0.02016607354685647  # F17
$ws.Cells.Item(17, 8).Value = 0.1553973902728351  # H17
$ws.Cells.Item(17, 9).Value = 0.08659549228944247  # I17
$ws.Cells.Item(17, 10).Value = 0.4163701067615658  # J17
$ws.Cells.Item(17, 11).Value = 0.1138790035587189  # K17
$ws.Cells.Item(17, 13).Value = 0.01897983392645314  # M17
$ws.Cells.Item(17, 15).Value = 0.05931198102016608  # O17
$ws.Cells.Item(17, 19).Value = 0.129300118623962  # S17

# Row 18 (Starting_State = 16)
$ws.Cells.Item(18, 6).Value = 0.02472527472527472  # F18
$ws.Cells.Item(18, 8).Value = 0.1785714285714286  # H18
$ws.Cells.Item(18, 9).Value = 0.09615384615384616  # I18
$ws.Cells.Item(18, 10).Value = 0.4230769230769231  # J18
$ws.Cells.Item(18, 11).Value = 0.1181318681318681  # K18
$ws.Cells.Item(18, 13).Value = 0.01098901098901099  # M18
$ws.Cells.Item(18, 15).Value = 0.04395604395604396  # O18
$ws.Cells.Item(18, 19).Value = 0.1043956043956044  # S18

# Row 19 (Starting_State = 17)
$ws.Cells.Item(19, 6).Value = 0.0279089738085015  # F19
$ws.Cells.Item(19, 8).Value = 0.2091026191498497  # H19
$ws.Cells.Item(19, 9).Value = 0.0734220695577501  # I19
$ws.Cells.Item(19, 10).Value = 0.3593817088879347  # J19
$ws.Cells.Item(19, 11).Value = 0.1240875912408759  # K19
$ws.Cells.Item(19, 13).Value = 0.0176041219407471  # M19
$ws.Cells.Item(19, 14).Value = 0.0008587376556462001  # N19
$ws.Cells.Item(19, 15).Value = 0.0790038643194504  # O19
$ws.Cells.Item(19, 19).Value = 0.1086303134392443  # S19
